$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 371, pushing the existing 371-393 block
# down to 374-396 (formats/styles of row 371 are inherited by the
# inserted rows, matching column D's date style).
$ws.Rows.Item(371).Resize(3).Insert()

# New row 371: Especial, 16 units
$ws.Cells.Item(371, 1).Value = 8
$ws.Cells.Item(371, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(371, 3).Value = "Coquimbo"
$ws.Cells.Item(371, 4).Value = 44746
$ws.Cells.Item(371, 5).Value = 4
$ws.Cells.Item(371, 6).Value = "Fruta"
$ws.Cells.Item(371, 7).Value = 100101
$ws.Cells.Item(371, 8).Value = "Berries"
$ws.Cells.Item(371, 9).Value = 100101007
$ws.Cells.Item(371, 10).Value = "Kiwi"
$ws.Cells.Item(371, 11).Value = "Hayward"
$ws.Cells.Item(371, 12).Value = "Especial"
$ws.Cells.Item(371, 13).Value = 16
$ws.Cells.Item(371, 14).Value = 240000
$ws.Cells.Item(371, 15).Value = 250000
$ws.Cells.Item(371, 16).Value = 245000
$ws.Cells.Item(371, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(371, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(371, 19).Value = 544
$ws.Cells.Item(371, 20).Value = 450

# New row 372: Primera, 20 units
$ws.Cells.Item(372, 1).Value = 8
$ws.Cells.Item(372, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44746
$ws.Cells.Item(372, 5).Value = 4
$ws.Cells.Item(372, 6).Value = "Fruta"
$ws.Cells.Item(372, 7).Value = 100101
$ws.Cells.Item(372, 8).Value = "Berries"
$ws.Cells.Item(372, 9).Value = 100101007
$ws.Cells.Item(372, 10).Value = "Kiwi"
$ws.Cells.Item(372, 11).Value = "Hayward"
$ws.Cells.Item(372, 12).Value = "Primera"
$ws.Cells.Item(372, 13).Value = 20
$ws.Cells.Item(372, 14).Value = 210000
$ws.Cells.Item(372, 15).Value = 220000
$ws.Cells.Item(372, 16).Value = 215000
$ws.Cells.Item(372, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(372, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(372, 19).Value = 478
$ws.Cells.Item(372, 20).Value = 450

# New row 373: Segunda, 20 units
$ws.Cells.Item(373, 1).Value = 8
$ws.Cells.Item(373, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(373, 3).Value = "Coquimbo"
$ws.Cells.Item(373, 4).Value = 44746
$ws.Cells.Item(373, 5).Value = 4
$ws.Cells.Item(373, 6).Value = "Fruta"
$ws.Cells.Item(373, 7).Value = 100101
$ws.Cells.Item(373, 8).Value = "Berries"
$ws.Cells.Item(373, 9).Value = 100101007
$ws.Cells.Item(373, 10).Value = "Kiwi"
$ws.Cells.Item(373, 11).Value = "Hayward"
$ws.Cells.Item(373, 12).Value = "Segunda"
$ws.Cells.Item(373, 13).Value = 20
$ws.Cells.Item(373, 14).Value = 160000
$ws.Cells.Item(373, 15).Value = 170000
$ws.Cells.Item(373, 16).Value = 165000
$ws.Cells.Item(373, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(373, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(373, 19).Value = 367
$ws.Cells.Item(373, 20).Value = 450
